# Burn down chart finished
#
# Sprint 2 actually wrapped up on day 18 instead of the originally planned
# day 20, so the tail of the "Sprint 2" data table needs correcting:
#   - C16:C19 ("Ideal tasks rem.") were driven by a shared formula; now that
#     the sprint is closed out they're fixed to their final literal values
#   - D18:D19 ("Actual tasks rem.") are corrected so actual work reaches 0
#     by day 18
#   - days 19 and 20 (rows 20 and 21) never happened, so their data is
#     cleared (the rows/date formatting stay in place)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")

$ws.Range("C16").Value = 7
$ws.Range("C17").Value = 5
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0

# Rows 20-21 (days 19-20) no longer have data - clear the cell contents but
# keep the rows (and column A's date-number style) in place.
$ws.Range("A20:C21").ClearContents()
